{"js": "// Split the three run-on \"1)...2)...3)...\" / reference-list paragraphs into\n// line-broken segments: <w:t>segment</w:t><w:br/><w:t>segment</w:t>... inside\n// a single run, by inserting the replacement text with an embedded vertical\n// tab (U+000B) between segments \u2014 that is how Word represents a manual line\n// break (<w:br/>) run-internally, and Range.insertText(\"Replace\") rewrites\n// the paragraph's run(s) to match it.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Index 13: Portuguese \"Programa\" numbered list.\nconst programaPt = [\n  \"1) Refra\u00e7\u00e3o e reflex\u00e3o.\",\n  \"2) Espelhos planos e esf\u00e9ricos e lentes delgadas.\",\n  \"3) Polariza\u00e7\u00e3o.\",\n  \"4) Interfer\u00eancia de ondas planas.\",\n  \"5) Difra\u00e7\u00e3o.\",\n  \"6) Espectroscopia \u00f3tica.\",\n  \"7) Determina\u00e7\u00e3o da constante de Planck.\",\n  \"8) Radia\u00e7\u00e3o de corpo negro.\",\n].join(\"\\v\");\n\n// Index 14: English (italic) \"Programa\" numbered list.\nconst programaEn = [\n  \"1) Refraction and reflection.\",\n  \"2) Mirrors plans and spherical and thin lenses.\",\n  \"3) Polarization.\",\n  \"4) Plane wave interference.\",\n  \"5) Diffraction.\",\n  \"6) Optical Spectroscopy.\",\n  \"7) Planck Constant Determination.\",\n  \"8) Black-body radiation.\",\n].join(\"\\v\");\n\n// Index 18: \"Bibliografia\" reference list.\nconst bibliografia = [\n  \"Apostilas do Laborat\u00f3rio de Ensino de F\u00edsica do IFSC/USP.\",\n  \"RESNICK, R.; HALLIDAY, D. Fundamentos de F\u00edsica. Vol. 4, LTC (2008).\",\n  \"TIPLER, P.; MOSCA, G. F\u00edsica para Cientistas e Engenheiros. Vol. 4, LTC (2008).\",\n  \"SEARS, F. W.; ZEMANSKY, M. W.; YOUNG, H. D.; FREEDMAN, R. A. F\u00edsica I, Vol. 4, Pearson Addison Wesley (2009).\",\n  \"JEWETT Jr, John W.; SERWAY, Raymond A. Princ\u00edpios de F\u00edsica. Vol. 4, Thomson Pioneira (2008).\",\n].join(\"\\v\");\n\nconst targets = {};\nfor (const p of paragraphs.items) {\n  p.load(\"text\");\n}\nawait context.sync();\n\nfor (const p of paragraphs.items) {\n  const text = p.text;\n  if (text.indexOf(\"1) Refra\u00e7\u00e3o e reflex\u00e3o.\") === 0) {\n    targets.pt = p;\n  } else if (text.indexOf(\"1) Refraction and reflection.\") === 0) {\n    targets.en = p;\n  } else if (text.indexOf(\"Apostilas do Laborat\u00f3rio de Ensino de F\u00edsica do IFSC/USP.\") === 0) {\n    targets.bib = p;\n  }\n}\n\nif (!targets.pt || !targets.en || !targets.bib) {\n  throw new Error(\"Could not locate one or more target paragraphs.\");\n}\n\ntargets.pt.getRange().insertText(programaPt, Word.InsertLocation.replace);\ntargets.en.getRange().insertText(programaEn, Word.InsertLocation.replace);\ntargets.bib.getRange().insertText(bibliografia, Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# Split the three run-on \"1)...2)...3)...\" / reference-list paragraphs into\n# line-broken segments. A backtick-v (`v) is PowerShell's escape for the\n# vertical-tab character (U+000B), which Word represents as a manual line\n# break (<w:br/>) when it appears inside Range.Text \u2014 assigning such a string\n# to a paragraph's Range.Text rewrites its run(s) into\n# <w:t>segment</w:t><w:br/><w:t>segment</w:t>... within a single run.\n\n$d = $word.ActiveDocument\n\n$programaPt = \"1) Refra\u00e7\u00e3o e reflex\u00e3o.`v2) Espelhos planos e esf\u00e9ricos e lentes delgadas.`v3) Polariza\u00e7\u00e3o.`v4) Interfer\u00eancia de ondas planas.`v5) Difra\u00e7\u00e3o.`v6) Espectroscopia \u00f3tica.`v7) Determina\u00e7\u00e3o da constante de Planck.`v8) Radia\u00e7\u00e3o de corpo negro.\"\n\n$programaEn = \"1) Refraction and reflection.`v2) Mirrors plans and spherical and thin lenses.`v3) Polarization.`v4) Plane wave interference.`v5) Diffraction.`v6) Optical Spectroscopy.`v7) Planck Constant Determination.`v8) Black-body radiation.\"\n\n$bibliografia = \"Apostilas do Laborat\u00f3rio de Ensino de F\u00edsica do IFSC/USP.`vRESNICK, R.; HALLIDAY, D. Fundamentos de F\u00edsica. Vol. 4, LTC (2008).`vTIPLER, P.; MOSCA, G. F\u00edsica para Cientistas e Engenheiros. Vol. 4, LTC (2008).`vSEARS, F. W.; ZEMANSKY, M. W.; YOUNG, H. D.; FREEDMAN, R. A. F\u00edsica I, Vol. 4, Pearson Addison Wesley (2009).`vJEWETT Jr, John W.; SERWAY, Raymond A. Princ\u00edpios de F\u00edsica. Vol. 4, Thomson Pioneira (2008).\"\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text\n    if ($t.StartsWith(\"1) Refra\u00e7\u00e3o e reflex\u00e3o.\")) {\n        $p.Range.Text = $programaPt\n    }\n    elseif ($t.StartsWith(\"1) Refraction and reflection.\")) {\n        $p.Range.Text = $programaEn\n    }\n    elseif ($t.StartsWith(\"Apostilas do Laborat\u00f3rio de Ensino de F\u00edsica do IFSC/USP.\")) {\n        $p.Range.Text = $bibliografia\n    }\n}\n"}
